# Update imputed values in the KNN result data sheet
# (commit message: 'Update Name of Algo' refers to re-running the KNN imputation
#  algorithm on column A/B/C, producing refreshed imputed values)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4" = -21.146
    "A6" = -22.283
    "A7" = -21.191
    "B7" = 6.765000000000001
    "A8" = -21.672
    "B11" = 6.734
    "B12" = 5.357
    "C12" = -11.692
    "C13" = -13.281
    "C14" = -12.161
    "B15" = 5.086999999999999
    "A16" = -21.949
    "C16" = -13.14
    "C19" = -12.635
    "A20" = -21.723
    "B20" = 5.916
    "C20" = -12.836
    "A21" = -21.536
    "B21" = 7.053
    "B22" = 8.516000000000002
    "C22" = -12.393
    "B23" = 7.533999999999999
    "A28" = -21.929
    "A29" = -21.698
    "B29" = 5.817
    "A30" = -21.807
    "A32" = -21.685
    "B34" = 8.300999999999998
    "C36" = -12.682
    "A40" = -20.664
    "B42" = 7.392
    "B43" = 6.282999999999999
    "C43" = -13.17
    "B44" = 5.993
    "B45" = 5.165999999999999
    "A46" = -21.108
    "B46" = 6.264999999999999
    "C46" = -13.338
    "B50" = 5.100000000000001
    "C50" = -13.419
    "A51" = -21.163
    "B51" = 7.514
    "A52" = -21.484
    "A57" = -21.432
    "B57" = 6.628
    "A59" = -22.415
    "A62" = -21.96
    "B65" = 5.306
    "A66" = -21.64
    "B66" = 5.473
    "B67" = 5.8
    "A73" = -20.731
    "A74" = -21.195
    "C76" = -11.854
    "A77" = -20.934
    "B79" = 5.53
    "B84" = 6.209
    "B87" = 4.447000000000001
    "A92" = -21.53
    "B92" = 5.527
    "C95" = -11.862
    "B97" = 7.121
    "C97" = -12.592
    "C99" = -12.123
    "A100" = -21.481
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

